$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "org-customprop-scanner"
$ws.Range("B2").Value = "https://github.com/Sohanuu66/org-customprop-scanner"
$ws.Range("C2").Value = $true
$ws.Range("D2").Value = "changes required"
